# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits the leading "word + trailing space" run in the title and
# caption textboxes of both slides back into two separate runs: the
# word run and a standalone " " run, leaving the following run(s)
# untouched.

function Split-LeadingWordRun($Shape, $Word) {
    $tr = $Shape.TextFrame.TextRange
    $len = $Word.Length
    $sub = $tr.Characters(1, $len)
    # Re-assigning the identical text forces the host to materialize the
    # selected character span as its own run, splitting it away from the
    # trailing " " (and whatever run(s) follow it).
    $sub.Text = $Word
}

$p = $ppt.ActivePresentation

$slide1 = $p.Slides.Item(1)
Split-LeadingWordRun $slide1.Shapes.Item(1) "Slide"
Split-LeadingWordRun $slide1.Shapes.Item(3) "an"

$slide2 = $p.Slides.Item(2)
Split-LeadingWordRun $slide2.Shapes.Item(1) "Slide"
Split-LeadingWordRun $slide2.Shapes.Item(4) "an"
